# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) / "Valor Mora" (column F) block (rows 16-61)
# is re-ordered from descending (2307 .. 1910) to ascending (1910 .. 2307).
# Net effect: column E gets the reversed list of period codes and the one
# non-standard "Valor Mora" amount (34666, all the rest being 40000) travels
# from the row that used to hold "2307" (now the earliest period, 1910) to
# the row that now holds "2307" (the latest period).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periods = @(
    "1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
    "2301","2302","2303","2304","2305","2306","2307"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    if ($periods[$i] -eq "2307") {
        $ws.Cells.Item($row, 6).Value = 34666
    } else {
        $ws.Cells.Item($row, 6).Value = 40000
    }
}
